# Updates the "cryptos" price/volume table (Sheet1) to the latest scraped
# values. Price cells (column D) that would otherwise be auto-parsed by
# Excel as numbers are written with a leading apostrophe so they stay text,
# matching the sheet's original inline-string ("t=inlineStr") formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.905.86"
$ws.Range("E2").Value = "  -3.42%  "
$ws.Range("D3").Value = "2.489.18"
$ws.Range("E3").Value = "  -5.98%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'557.89"
$ws.Range("E5").Value = "  -4.20%  "
$ws.Range("D6").Value = "'148.62"
$ws.Range("E6").Value = "  -5.14%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.603"
$ws.Range("E8").Value = "  -4.01%  "
$ws.Range("D9").Value = "2.488.61"
$ws.Range("E9").Value = "  -5.92%  "
$ws.Range("E10").Value = "  -7.99%  "
$ws.Range("D11").Value = "'5.52"
$ws.Range("E11").Value = "  -5.30%  "
$ws.Range("E12").Value = "  -1.40%  "
$ws.Range("D13").Value = "'0.361"
$ws.Range("E13").Value = "  -6.38%  "
$ws.Range("D14").Value = "'26.70"
$ws.Range("E14").Value = "  -6.86%  "
$ws.Range("D15").Value = "2.933.98"
$ws.Range("E15").Value = "  -6.11%  "
$ws.Range("D16").Value = "'0.0000170"
$ws.Range("E16").Value = "  -8.36%  "
$ws.Range("D17").Value = "61.753.09"
$ws.Range("E17").Value = "  -3.50%  "
$ws.Range("D18").Value = "2.486.13"
$ws.Range("E18").Value = "  -6.26%  "
$ws.Range("D19").Value = "'11.29"
$ws.Range("E19").Value = "  -7.71%  "
$ws.Range("D20").Value = "'7.22"
$ws.Range("E20").Value = "  -7.08%  "
$ws.Range("D21").Value = "'4.26"
$ws.Range("E21").Value = "  -6.46%  "
$ws.Range("D22").Value = "'324.26"
$ws.Range("E22").Value = "  -6.51%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "'1.91"
$ws.Range("E24").Value = "  +2.47%  "
$ws.Range("D25").Value = "'64.67"
$ws.Range("E25").Value = "  -5.27%  "
$ws.Range("D26").Value = "'0.0000103"
$ws.Range("E26").Value = "  -8.85%  "
$ws.Range("D27").Value = "'576.48"
$ws.Range("E27").Value = "  -2.10%  "
$ws.Range("D28").Value = "2.604.46"
$ws.Range("E28").Value = "  -6.23%  "
$ws.Range("E29").Value = "  -6.33%  "
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'8.40"
$ws.Range("E31").Value = "  -10.02%  "
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").Value = "'7.87"
$ws.Range("E32").Value = "  -4.84%  "
$ws.Range("D33").Value = "'0.152"
$ws.Range("E33").Value = "  -6.47%  "
$ws.Range("E34").Value = "  -6.63%  "
$ws.Range("D35").Value = "'1.61"
$ws.Range("E35").Value = "  -7.99%  "
$ws.Range("D36").Value = "'6.03"
$ws.Range("E36").Value = "  -9.58%  "
$ws.Range("D37").Value = "'4.99"
$ws.Range("E37").Value = "  -9.77%  "
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("D39").Value = "'0.386"
$ws.Range("E39").Value = "  -4.55%  "
$ws.Range("D40").Value = "'18.71"
$ws.Range("E40").Value = "  -5.60%  "
$ws.Range("D41").Value = "'145.20"
$ws.Range("E41").Value = "  -4.30%  "
$ws.Range("D42").Value = "'1.79"
$ws.Range("E42").Value = "  -7.36%  "
$ws.Range("D44").Value = "'2.48"
$ws.Range("E44").Value = "  -4.05%  "
$ws.Range("D45").Value = "'40.71"
$ws.Range("E45").Value = "  -2.89%  "
$ws.Range("D46").Value = "'149.88"
$ws.Range("E46").Value = "  -8.64%  "
$ws.Range("D47").Value = "'3.68"
$ws.Range("E47").Value = "  -6.33%  "
$ws.Range("D48").Value = "'22.15"
$ws.Range("E48").Value = "  -9.70%  "
$ws.Range("D49").Value = "'0.0546"
$ws.Range("E49").Value = "  -7.86%  "
$ws.Range("D50").Value = "'0.601"
$ws.Range("E50").Value = "  -5.55%  "
$ws.Range("D51").Value = "'0.0949"
$ws.Range("E51").Value = "  -5.55%  "
